{"js": "// The document has a red/bold \"features\" sub-heading (right before the\n// \"Achievements\" line) that should read \"Features\" (capitalized first\n// letter), matching the rest of the headings (\"Controls\", \"Install\n// Instructions\", \"Known Problems\") which are all capitalized.\nconst body = context.document.body;\n\n// matchWholeWord + matchCase ensures we only touch the standalone heading\n// word \"features\" and not some other occurrence elsewhere in the document.\nconst results = body.search(\"features\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  // Replacing in place preserves the run's existing formatting\n  // (bold, red color, size 32).\n  results.items[i].insertText(\"Features\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document has a red/bold \"features\" sub-heading (right before the\n# \"Achievements\" line) that should read \"Features\" (capitalized first\n# letter), matching the rest of the headings (\"Controls\", \"Install\n# Instructions\", \"Known Problems\") which are all capitalized.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n# MatchWholeWord + MatchCase ensures we only touch the standalone heading\n# word \"features\" and not some other occurrence elsewhere in the document.\n# Replacing via Find/Execute keeps the existing run formatting (bold, red,\n# size 32) on the matched text.\n$result = $find.Execute(\n    \"features\",   # FindText\n    $true,        # MatchCase\n    $true,        # MatchWholeWord\n    $false,       # MatchWildcards\n    $false,       # MatchSoundsLike\n    $false,       # MatchAllWordForms\n    $true,        # Forward\n    1,            # Wrap (wdFindContinue)\n    $false,       # Format\n    \"Features\",   # ReplaceWith\n    2             # Replace (wdReplaceAll)\n)\n"}
